# "Update for new pcb"
# Re-labels the microcontroller pins (new PCB uses different D-pin / A-pin
# naming) and adds a small "Verfügbare Pins:" (available pins) legend table
# in column I so the new pin names used in column B are documented.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) New "Verfügbare Pins:" legend table in column I (filled in first, top
#    to bottom, so the shared-string table gets the same ordering as the
#    authored workbook).
# ---------------------------------------------------------------------
$ws.Range("I1").Value = "Verfügbare Pins:"

$ws.Range("I3").Value = "D4"
$ws.Range("I4").Value = "D3"
$ws.Range("I5").Value = "D2"
$ws.Range("I6").Value = "A3"
$ws.Range("I7").Value = "A2"
$ws.Range("I8").Value = "A1"
$ws.Range("I9").Value = "A0"
$ws.Range("I10").Value = "D5"
$ws.Range("I11").Value = "D6"
$ws.Range("I12").Value = "D7"
$ws.Range("I13").Value = "D8"
$ws.Range("I14").Value = "D9"
$ws.Range("I15").Value = "D10"
$ws.Range("I16").Value = "D11"

# Header cell uses the same style as the other "Legende:" header (G1):
# centered both horizontally and vertically (no fill).
$ws.Range("G1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "Verfügbare Pins:"

# Plain centered cells (no fill) for I3, I5, I6, I7, I8, I9, I12, I13
$ws.Range("I3").HorizontalAlignment = -4108
$ws.Range("I5").HorizontalAlignment = -4108
$ws.Range("I6").HorizontalAlignment = -4108
$ws.Range("I7").HorizontalAlignment = -4108
$ws.Range("I8").HorizontalAlignment = -4108
$ws.Range("I9").HorizontalAlignment = -4108
$ws.Range("I12").HorizontalAlignment = -4108
$ws.Range("I13").HorizontalAlignment = -4108

# Highlighted (filled + centered) cells for I4, I10, I11, I14, I15, I16 --
# copy the existing highlight fill used elsewhere on the sheet (G3) and
# then center the text.
$fillSource = $ws.Range("G3")
$fillSource.Copy()

$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("I10").PasteSpecial(-4122)
$ws.Range("I11").PasteSpecial(-4122)
$ws.Range("I14").PasteSpecial(-4122)
$ws.Range("I15").PasteSpecial(-4122)
$ws.Range("I16").PasteSpecial(-4122)

$ws.Range("I4").HorizontalAlignment = -4108
$ws.Range("I10").HorizontalAlignment = -4108
$ws.Range("I11").HorizontalAlignment = -4108
$ws.Range("I14").HorizontalAlignment = -4108
$ws.Range("I15").HorizontalAlignment = -4108
$ws.Range("I16").HorizontalAlignment = -4108

# Column I width, wide enough to fit "Verfügbare Pins:".
$ws.Columns.Item(9).ColumnWidth = 14.86

# ---------------------------------------------------------------------
# 2) Relabel the "Pin" column (B) for every connector block (X8 / I2C-Addr
#    8, 9, 10, 11) using the new pin names.
# ---------------------------------------------------------------------

# Block 8 (rows 3-16)
$ws.Range("B3").Value = "D3"
$ws.Range("B4").Value = "D4"
$ws.Range("B5").Value = "D5"
$ws.Range("B6").Value = "D6"
$ws.Range("B7").Value = "D2"
$ws.Range("B8").Value = "D9"
$ws.Range("B9").Value = "D10"
$ws.Range("B10").Value = "D11"
$ws.Range("B11").Value = "A3"
$ws.Range("B12").Value = "A2"
$ws.Range("B13").Value = "A1"
$ws.Range("B14").Value = "A0"
$ws.Range("B15").Value = "D7"
$ws.Range("B16").Value = "D8"

# Block 9 (rows 18-31)
$ws.Range("B18").Value = "D3"
$ws.Range("B19").Value = "D4"
$ws.Range("B20").Value = "D5"
$ws.Range("B21").Value = "D6"
$ws.Range("B22").Value = "D2"
$ws.Range("B23").Value = "D9"
$ws.Range("B24").Value = "D10"
$ws.Range("B25").Value = "D11"
$ws.Range("B26").Value = "A3"
$ws.Range("B27").Value = "A2"
$ws.Range("B28").Value = "A1"
$ws.Range("B29").Value = "A0"
$ws.Range("B30").Value = "D7"
$ws.Range("B31").Value = "D8"

# Block 10 (rows 33-46)
$ws.Range("B33").Value = "D3"
$ws.Range("B34").Value = "D4"
$ws.Range("B35").Value = "D5"
$ws.Range("B36").Value = "D6"
$ws.Range("B37").Value = "D2"
$ws.Range("B38").Value = "D9"
$ws.Range("B39").Value = "D10"
$ws.Range("B40").Value = "D11"
$ws.Range("B41").Value = "A3"
$ws.Range("B42").Value = "A2"
$ws.Range("B43").Value = "A1"
$ws.Range("B44").Value = "A0"
$ws.Range("B45").Value = "D7"
$ws.Range("B46").Value = "D8"

# Block 11 (rows 48-61)
$ws.Range("B48").Value = "D3"
$ws.Range("B49").Value = "D4"
$ws.Range("B50").Value = "D5"
$ws.Range("B51").Value = "D6"
$ws.Range("B52").Value = "D2"
$ws.Range("B53").Value = "D9"
$ws.Range("B54").Value = "D10"
$ws.Range("B55").Value = "D11"
$ws.Range("B56").Value = "A3"
$ws.Range("B57").Value = "A2"
$ws.Range("B58").Value = "A1"
$ws.Range("B59").Value = "A0"
$ws.Range("B60").Value = "D7"
$ws.Range("B61").Value = "D8"

# ---------------------------------------------------------------------
# 3) Restore the view/selection state recorded in the authored workbook.
# ---------------------------------------------------------------------
$ws.Range("G32").Select()
